# Reorder the "Recorded By" (column G) author lists so that a rotated
# (non-"System"-first) list gets its first entry moved to the end.
#
# For every data row (2..157) in column G, if the cell contains a
# comma-separated list with more than one entry and the first entry is not
# exactly "System", move that first entry to the end of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ",\s*"
    if ($parts.Count -le 1) { continue }
    if ($parts[0] -eq "System") { continue }

    $newParts = $parts[1..($parts.Count - 1)] + $parts[0]
    $newVal = [string]::Join(", ", $newParts)

    $cell.Value2 = $newVal
}
